$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Fri Mar 08 00:39:27 EST 2024"
$ws.Range("B3").Value = "Fri Mar 08 00:40:05 EST 2024"
$ws.Range("B4").Value = "Fri Mar 08 00:40:29 EST 2024"
$ws.Range("B5").Value = "Fri Mar 08 00:41:07 EST 2024"
$ws.Range("B6").Value = "Fri Mar 08 00:41:44 EST 2024"
$ws.Range("B7").Value = "Fri Mar 08 00:42:14 EST 2024"
$ws.Range("B8").Value = "Fri Mar 08 00:43:08 EST 2024"
$ws.Range("B9").Value = "Fri Mar 08 00:43:41 EST 2024"
$ws.Range("B10").Value = "Fri Mar 08 00:44:15 EST 2024"
$ws.Range("B11").Value = "Fri Mar 08 00:44:41 EST 2024"
$ws.Range("B12").Value = "Fri Mar 08 00:45:21 EST 2024"
$ws.Range("B13").Value = "Fri Mar 08 00:46:02 EST 2024"
$ws.Range("B14").Value = "Fri Mar 08 00:46:52 EST 2024"
$ws.Range("B15").Value = "Fri Mar 08 00:47:23 EST 2024"
$ws.Range("B16").Value = "Fri Mar 08 00:48:06 EST 2024"
$ws.Range("B17").Value = "Fri Mar 08 00:48:42 EST 2024"
$ws.Range("B18").Value = "Fri Mar 08 00:49:20 EST 2024"
$ws.Range("B19").Value = "Fri Mar 08 00:50:00 EST 2024"
$ws.Range("B20").Value = "Fri Mar 08 00:50:37 EST 2024"
$ws.Range("B28").Value = "Fri Mar 08 00:51:07 EST 2024"
$ws.Range("B29").Value = "Fri Mar 08 00:51:36 EST 2024"
$ws.Range("B30").Value = "Fri Mar 08 00:52:07 EST 2024"
$ws.Range("B31").Value = "Fri Mar 08 00:52:41 EST 2024"
$ws.Range("B32").Value = "Fri Mar 08 00:53:25 EST 2024"
$ws.Range("B33").Value = "Fri Mar 08 00:54:10 EST 2024"
$ws.Range("B34").Value = "Fri Mar 08 00:54:57 EST 2024"
$ws.Range("B35").Value = "Fri Mar 08 00:55:53 EST 2024"
$ws.Range("B36").Value = "Fri Mar 08 00:56:42 EST 2024"
$ws.Range("B37").Value = "Fri Mar 08 00:57:29 EST 2024"
$ws.Range("B38").Value = "Fri Mar 08 00:58:13 EST 2024"
$ws.Range("B39").Value = "Fri Mar 08 00:58:55 EST 2024"
$ws.Range("B40").Value = "Fri Mar 08 00:59:38 EST 2024"
$ws.Range("B41").Value = "Fri Mar 08 01:00:15 EST 2024"
$ws.Range("B42").Value = "Fri Mar 08 01:00:52 EST 2024"
$ws.Range("B43").Value = "Fri Mar 08 01:01:15 EST 2024"
$ws.Range("B44").Value = "Fri Mar 08 01:02:08 EST 2024"
$ws.Range("B45").Value = "Fri Mar 08 01:02:57 EST 2024"
$ws.Range("B46").Value = "Fri Mar 08 01:03:33 EST 2024"
$ws.Range("B47").Value = "Fri Mar 08 01:04:00 EST 2024"
$ws.Range("B48").Value = "Fri Mar 08 01:04:40 EST 2024"
$ws.Range("B49").Value = "Fri Mar 08 01:05:18 EST 2024"
$ws.Range("B50").Value = "Fri Mar 08 01:05:55 EST 2024"
$ws.Range("B51").Value = "Fri Mar 08 01:06:38 EST 2024"
$ws.Range("B52").Value = "Fri Mar 08 01:07:01 EST 2024"
$ws.Range("B53").Value = "Fri Mar 08 01:07:40 EST 2024"
$ws.Range("B54").Value = "Fri Mar 08 01:08:13 EST 2024"
